$d = $word.ActiveDocument

# --- Locate the target paragraph -------------------------------------------------
# This is the last bullet in the "Notes" list, currently reading:
#   "Weapons/Items als ScriptableObjects, dann sind die auch verwendet"
# and carrying the (hidden) "_GoBack" bookmark.
$targetPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*als ScriptableObjects, dann sind die auch verwendet*") {
        $targetPara = $para
    }
}

# --- Split off a brand-new list paragraph after it --------------------------------
# InsertParagraphAfter() on a collapsed range at the end of the paragraph creates a
# sibling <w:p> that inherits the same pPr (Listenabsatz / numPr ilvl0 numId2).
$tailRange = $targetPara.Range
$tailRange.Collapse(0)
$tailRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
$newRange.Collapse(0)
$newRange.InsertBefore("Fighter nicht tot sondern lastBreath wenn Health = 0 gefangen durch Tests")

# --- Temporarily pad the document -------------------------------------------------
# Placing a new bookmark at/within the last two character positions of the whole
# document confuses this host's bookmark resolver (it silently resets to (0,0)-ish).
# Append two throw-away characters after our new text so the insertion point we need
# is no longer at the very tail of the document, then strip them back out afterwards.
$padRange = $d.Paragraphs.Last.Range
$padRange.Collapse(0)
$padRange.InsertAfter("ZZ")

# --- Move the "_GoBack" bookmark onto the end of the new paragraph ----------------
$insertionPoint = $d.Paragraphs.Last.Range.End - 1 - 2   # -1 skips the pilcrow, -2 skips "ZZ"
$bookmarkRange = $d.Range($insertionPoint, $insertionPoint)

$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# --- Remove the padding characters again -------------------------------------------
$finalParaEnd = $d.Paragraphs.Last.Range.End
$padRange2 = $d.Range($finalParaEnd - 1 - 2, $finalParaEnd - 1)
$padRange2.Delete()
